$d = $word.ActiveDocument

$replacements = @(
    @("113÷9=", "970÷7="),
    @("344÷6=", "535÷6="),
    @("724÷7=", "432÷8="),
    @("960÷3=", "256÷8="),
    @("204÷8=", "984÷5="),
    @("823÷2=", "675÷8="),
    @("216÷2=", "328÷2="),
    @("342÷8=", "346÷6="),
    @("809÷5=", "196÷8="),
    @("639÷7=", "917÷9="),
    @("556÷4=", "394÷9="),
    @("903÷2=", "524÷8="),
    @("394÷8=", "266÷3="),
    @("410÷3=", "628÷4="),
    @("506÷8=", "783÷4="),
    @("163÷4=", "125÷9="),
    @("452÷4=", "739÷2="),
    @("745÷5=", "670÷6="),
    @("531÷2=", "917÷4="),
    @("484÷3=", "708÷8="),
    @("563÷2=", "969÷3="),
    @("267÷8=", "497÷6="),
    @("214÷6=", "701÷2="),
    @("168÷9=", "212÷2="),
    @("748÷2=", "857÷6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
